$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> (nombre_aides, montant_total)
$updates = @{
    35  = @(19274, 61939679)
    91  = @(151137, 482422118)
    92  = @(409117, 1595230363)
    93  = @(209581, 1308919641)
    95  = @(50766, 932507917)
    96  = @(17277, 793121202)
    104 = @(135240, 272208155)
    119 = @(351, 10119042)
    145 = @(11832, 182727133)
    164 = @(34185, 117606827)
    183 = @(34, 5661204)
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Cells.Item($row, 3).Value = $vals[0]
    $ws.Cells.Item($row, 5).Value = $vals[1]
}
